$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "... password protecte" + _GoBack bookmark + "d." runs back
#    into a single run/text "... password protected." This also removes the
#    old _GoBack bookmark that previously sat at this location.
# ---------------------------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute("password protected.", $false, $false, $false, $false, $false, $true, 1, $false, "password protected.", 2)

# ---------------------------------------------------------------------------
# 2) Re-insert the _GoBack bookmark at its new location, right after
#    "...it would be beneficial for" (splitting that run in two). Since
#    bookmark names must be unique, adding it here moves it from wherever it
#    was (it no longer exists after step 1) to this new spot.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
[void]$rng2.Find.Execute("it would be beneficial for", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($rng2.End, $rng2.End)
[void]$d.Bookmarks.Add("_GoBack", $splitPoint)

# ---------------------------------------------------------------------------
# 3) Collapse the five PAGEREF TOC fields whose target page changed into
#    plain text runs carrying the updated page number (the bookmarks/
#    headings they reference are unaffected).
# ---------------------------------------------------------------------------
$tocUpdates = @{
    "_Toc18551422" = "2"
    "_Toc18551427" = "3"
    "_Toc18551428" = "3"
    "_Toc18551429" = "3"
    "_Toc18551432" = "4"
}

for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $f = $d.Fields.Item($i)
    $code = $f.Code.Text
    foreach ($key in $tocUpdates.Keys) {
        if ($code -like "*$key*") {
            $newPage = $tocUpdates[$key]
            $resultRange = $f.Result
            [void]$resultRange.Find.Execute($resultRange.Text, $false, $false, $false, $false, $false, $true, 1, $false, $newPage, 2)
            $f.Unlink()
        }
    }
}

# ---------------------------------------------------------------------------
# 4) Update the cached value of the header's PAGE field from 5 to 2 (field
#    code stays intact, only the displayed/cached result text changes).
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRange = $hdr.Range
for ($i = 1; $i -le $hdrRange.Fields.Count; $i++) {
    $hf = $hdrRange.Fields.Item($i)
    if ($hf.Code.Text -like "*PAGE*") {
        $hres = $hf.Result
        [void]$hres.Find.Execute($hres.Text, $false, $false, $false, $false, $false, $true, 1, $false, "2", 2)
    }
}
